# add : switch show tree on tag item page
# Append two new "add" log rows for the new `tag` entity (images_aeriennes_1
# and images_aeriennes_2), following the existing log-table layout:
#   A=timestamp  B=type  C=entity  D=entity_id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$newRow1 = $lastRow + 1
$ws.Cells.Item($newRow1, 1).Value = 1749897806
$ws.Cells.Item($newRow1, 2).Value = "add"
$ws.Cells.Item($newRow1, 3).Value = "tag"
$ws.Cells.Item($newRow1, 4).Value = "images_aeriennes_1"

$newRow2 = $newRow1 + 1
$ws.Cells.Item($newRow2, 1).Value = 1749897806
$ws.Cells.Item($newRow2, 2).Value = "add"
$ws.Cells.Item($newRow2, 3).Value = "tag"
$ws.Cells.Item($newRow2, 4).Value = "images_aeriennes_2"
